$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style (from an existing date cell) into the new
# A16:A20 cells first, so the number format/border/font match the rest of
# column A without creating a brand-new style entry.
$ws.Range("A15").Copy($ws.Range("A16:A20"))

# Row 16: Algyo / House
$ws.Range("A16").Value = 43847
$ws.Range("B16").Value = "Algyo"
$ws.Range("C16").Value = "House"
$ws.Range("D16").Value = 119.1052631578947
$ws.Range("E16").Value = 32676315.78947368
$ws.Range("F16").Value = 296798.3356605554
$ws.Range("G16").Value = 486.8421052631579
$ws.Range("H16").Value = 19

# Row 17: Morahalom / House
$ws.Range("A17").Value = 43847
$ws.Range("B17").Value = "Morahalom"
$ws.Range("C17").Value = "House"
$ws.Range("D17").Value = 115.8809523809524
$ws.Range("E17").Value = 25552380.95238095
$ws.Range("F17").Value = 228041.7572974148
$ws.Range("G17").Value = 181.5
$ws.Range("H17").Value = 42

# Row 18: Szeged / Flat (new category)
$ws.Range("A18").Value = 43847
$ws.Range("B18").Value = "Szeged"
$ws.Range("C18").Value = "Flat"
$ws.Range("D18").Value = 65.3952380952381
$ws.Range("E18").Value = 27610738.0952381
$ws.Range("F18").Value = 426755.1999621868
$ws.Range("H18").Value = 2520

# Row 19: Szeged / Garage
$ws.Range("A19").Value = 43847
$ws.Range("B19").Value = "Szeged"
$ws.Range("C19").Value = "Garage"
$ws.Range("D19").Value = 17.75
$ws.Range("E19").Value = 4148100
$ws.Range("F19").Value = 240649.1567228772
$ws.Range("H19").Value = 100

# Row 20: Szeged / House
$ws.Range("A20").Value = 43847
$ws.Range("B20").Value = "Szeged"
$ws.Range("C20").Value = "House"
$ws.Range("D20").Value = 163.6563658838072
$ws.Range("E20").Value = 56148529.04820766
$ws.Range("F20").Value = 979818.0885652011
$ws.Range("G20").Value = 460.5908529048208
$ws.Range("H20").Value = 809
